$wb = $excel.ActiveWorkbook

# --- Rename worksheets (sheet task order identifiers refreshed) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16509961825567062"
$wb.Worksheets.Item(2).Name = "NB_TO-16509961848607035"
$wb.Worksheets.Item(3).Name = "RS_TO-16509961848607035"
$wb.Worksheets.Item(4).Name = "TOL_TO-16509961849247053"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16509961849967093"

# --- Sheet 1 (GNG_TO) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16509961825167046.csv"
$ws1.Range("B3").Value = "GNG_stims-16509961825407057.csv"
$ws1.Range("B4").Value = "go_stims-16509961825407057.csv"
$ws1.Range("B5").Value = "GNG_stims-16509961825567062.csv"

# --- Sheet 2 (NB_TO) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16509961836287057.csv"
$ws2.Range("B3").Value = "ZB-match_5-16509961833967407.csv"
$ws2.Range("B4").Value = "TB-1650996184836705.csv"
$ws2.Range("B5").Value = "ZB-match_2-16509961828607552.csv"
$ws2.Range("B6").Value = "TB-16509961842847033.csv"
$ws2.Range("B7").Value = "TB-16509961847007456.csv"
$ws2.Range("B8").Value = "OB-16509961840527415.csv"
$ws2.Range("B9").Value = "OB-16509961835327418.csv"
$ws2.Range("B10").Value = "ZB-match_2-16509961832607038.csv"

# --- Sheet 3 (RS_TO) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4 (TOL_TO) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-1650996184892739.csv"
$ws4.Range("B3").Value = "ZM_stims-16509961848687057.csv"
$ws4.Range("B4").Value = "MM_stims-16509961849087396.csv"
$ws4.Range("B5").Value = "ZM_stims-1650996184892739.csv"
$ws4.Range("B6").Value = "MM_stims-16509961849247053.csv"
$ws4.Range("B7").Value = "ZM_stims-16509961849087396.csv"

# --- Sheet 5 (vSAT_TO) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16509961849407053.csv"
$ws5.Range("B3").Value = "SAT_stims-16509961849247053.csv"
$ws5.Range("B4").Value = "vSAT_stims-16509961849647412.csv"
$ws5.Range("B5").Value = "vSAT_stims-16509961849807403.csv"
